# Introduces the 'INTEGER' 'Cell' type.
#
# Row 1 of the sheet is a catalogue of example `Cell`s (one pair of
# columns per type: the example value, then a formula that concatenates
# it with itself). This change inserts a new "Integer" example (columns
# G/H: value 1 formatted as Text, formula =CONCAT(G1,G1)) right after the
# existing "Double" example (which slides from F to I), and adds a
# second "Double"-with-literal example (columns K/L: the same double
# value formatted as Text, formula =CONCAT(1,",1")). Everything that
# used to follow the Double example (Date, Currency, String, String
# w/red font, String w/alt font) shifts further to the right to make
# room, and the comment that was anchored on the first "String" example
# moves along with it (K1 -> O1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Phase 1: relocate the cells that slide right, right-to-left so a cell is
# always copied before anything overwrites it. .Copy carries the value,
# formula and formatting of each cell to its new home in one shot.
# ---------------------------------------------------------------------------

$ws.Range("M1").Copy($ws.Range("Q1"))   # "Automatic" example -> Q1
$ws.Range("L1").Copy($ws.Range("P1"))   # "String" (red font) example -> P1
$ws.Range("K1").Copy($ws.Range("O1"))   # "String" example -> O1
$ws.Range("H1").Copy($ws.Range("N1"))   # Currency example -> N1
$ws.Range("G1").Copy($ws.Range("M1"))   # Date example -> M1
$ws.Range("F1").Copy($ws.Range("I1"))   # Double example value -> I1
$excel.CutCopyMode = $false

# Move the comment that lived on the original "String" example (K1) along
# with its cell, to its new home at O1.
$ws.Range("K1").Comment.Delete()
$ws.Range("O1").AddComment("Note")

# ---------------------------------------------------------------------------
# Phase 2: write the new/changed content at the vacated columns. G1/H1/K1/
# L1 still carry left-over formatting from the cells that used to live
# there (Date/Currency/String), so their format is reset to plain/General
# first by pasting formats from a cell that has always been plain (E1).
# ---------------------------------------------------------------------------

$ws.Range("E1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# F1 now recomputes the relocated Double example's concatenation (the
# same formula the old "I1" used to carry for the first, Integer,
# example pair).
$ws.Range("F1").Formula = "=E1+E1"

# G1/H1: brand-new "Integer" example pair.
$ws.Range("G1").Value = 1
$ws.Range("G1").NumberFormat = "@"
$ws.Range("H1").Formula = "=CONCAT(G1,G1)"

# J1: concatenation sibling for the relocated Double example value (I1).
$ws.Range("J1").Formula = "=I1+I1"

# K1/L1: new second "Double" example pair, built from literal numbers.
$ws.Range("K1").Value = $ws.Range("I1").Value2
$ws.Range("K1").NumberFormat = "@"
$ws.Range("L1").Formula = '=CONCAT(1,",1")'

# The relocated Date cell (M1) already carries its m/d/yyyy format via
# the .Copy above. The relocated Currency cell (N1) picks up the new
# red-negative format code.
$ws.Range("N1").NumberFormat = '"R$"\ #,##0.00;[Red]\-"R$"\ #,##0.00'

# ---------------------------------------------------------------------------
# Column widths: the three narrow "value" columns of the Integer example
# (E/F/G) keep the original width, and every column from H onward is
# re-fit for its (shifted) content.
# ---------------------------------------------------------------------------
$ws.Columns("E:G").ColumnWidth = 1.1666666666666665
$ws.Columns("H:H").ColumnWidth = 2.1666666666666665
$ws.Columns("I:J").ColumnWidth = 3.1666666666666665
$ws.Columns("K:K").ColumnWidth = 3.0221354166666665
$ws.Columns("L:L").ColumnWidth = 2.7369791666666665
$ws.Columns("M:M").ColumnWidth = 9.877604166666666
$ws.Columns("N:N").ColumnWidth = 6.307291666666667
$ws.Columns("O:P").ColumnWidth = 5.307291666666667
$ws.Columns("Q:Q").ColumnWidth = 9.307291666666666

# ---------------------------------------------------------------------------
# Sheet-level bookkeeping to match the saved state: selection sits on the
# new Currency cell (N1).
# ---------------------------------------------------------------------------
$ws.Range("N1").Select()
